$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EBC")

# Add a new column J containing the R-R interval (column B) truncated to 2 decimal
# places. This is used for the 10- and 16-interval continuous-monitoring HRV
# "Estimated Breath Cycle" chart data field.
for ($r = 2; $r -le 17; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 10).Value = $excel.WorksheetFunction.Trunc($b, 2)
}

# Point the 10-interval Min/Max formulas at the new truncated column instead of
# the raw (untruncated) R-R interval column.
$ws.Range("C2").Formula = "=MIN(J2:J11)"
$ws.Range("D2").Formula = "=MAX(J2:J11)"

# Point the 16-interval Min formula at the new truncated column as well.
$ws.Range("F2").Formula = "=MIN(J2:J17)"

$ws.Calculate() | Out-Null

# Match the column width Excel settled on for the new data column.
$ws.Columns.Item(10).ColumnWidth = 21

$ws.Activate() | Out-Null
$ws.Range("E2").Select() | Out-Null
